$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.368.79"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.077.53"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +4.47%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.84"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.88%  "
$ws.Range("E6").Value = "  +2.13%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.28"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +4.58%  "
$ws.Range("E9").Value = "  +2.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.03"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.69%  "
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.383.03"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +4.61%  "
$ws.Range("E14").Value = "  +3.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.02"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("E16").Value = "  +2.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.24"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +4.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.061.22"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37.512.59"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.01"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +20.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.30"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "223.54"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.80%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.43"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.85%  "
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.80"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.86"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.80%  "
$ws.Range("E29").Value = "  +7.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.31"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("E31").Value = "  +5.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.119"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.47"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.35%  "
$ws.Range("E34").Value = "  +2.29%  "
$ws.Range("E35").Value = "  +10.05%  "
$ws.Range("E36").Value = "  +4.50%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  +14.78%  "
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("E40").Value = "  +0.49%  "
$ws.Range("E41").Value = "  -4.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0962"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +10.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.473.72"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.99%  "
$ws.Range("E44").Value = "  +19.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "94.86"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +7.76%  "
$ws.Range("E46").Value = "  +3.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.21"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +7.18%  "
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("E49").Value = "  +3.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.30"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +10.07%  "
$ws.Range("E51").Value = "  +1.80%  "
